$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected roster table (Player, Position, Team) for rows 2-19, derived from the
# target sharedStrings order + cell index mapping described by the diff.
$data = @(
    @("Fred VanVleet",          "PG",       "Houston Rockets"),
    @("Anfernee Simons",        "PG,SG",    "Portland Trail Blazers"),
    @("James Harden",           "PG,SG",    "LA Clippers"),
    @("Anthony Edwards",        "SG,SF",    "Minnesota Timberwolves"),
    @("Bilal Coulibaly",        "SG,SF",    "Washington Wizards"),
    @("Paul George",            "SG,SF,PF", "Philadelphia 76ers"),
    @("Giannis Antetokounmpo",  "PF,C",     "Milwaukee Bucks"),
    @("Wendell Carter Jr.",     "PF,C",     "Orlando Magic"),
    @("Nicolas Claxton",        "C",        "Brooklyn Nets"),
    @("Ivica Zubac",            "C",        "LA Clippers"),
    @("Jayson Tatum",           "SF,PF",    "Boston Celtics"),
    @("Keyonte George",         "PG,SG",    "Utah Jazz"),
    @("Amen Thompson",          "SG,SF",    "Houston Rockets"),
    @("Jaren Jackson Jr.",      "PF,C",     "Memphis Grizzlies"),
    @("Zion Williamson",        "PF,C",     "New Orleans Pelicans"),
    @("Jonathan Kuminga",       "SF,PF",    "Golden State Warriors"),
    @("Goga Bitadze",           "C",        "Orlando Magic"),
    @("Bradley Beal",           "PG,SG,SF", "Phoenix Suns")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
